$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mini")

# Update mark-level values (columns C, D, E for rows 4-10)
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = 80

$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 80

$ws.Range("D6").Value = 60
$ws.Range("E6").Value = 80

$ws.Range("C7").Value = 2000
$ws.Range("D7").Value = 3750
$ws.Range("E7").Value = 5000

$ws.Range("E8").Value = 65

$ws.Range("C9").Value = 500
$ws.Range("D9").Value = 1200
$ws.Range("E9").Value = 2000

$ws.Range("C10").Value = 5000
$ws.Range("D10").Value = 8000
$ws.Range("E10").Value = 10000

# Update the active cell selection to D8
$ws.Range("D8").Select()
